$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(77, 1).Value = "0083"
$ws.Cells.Item(77, 2).Value = "EASY"
$ws.Cells.Item(77, 3).Value = "Remove Duplicates from Sorted List"
$ws.Cells.Item(77, 4).Value = "method1.cpp"
$ws.Cells.Item(77, 5).Value = "Linked List"
$ws.Cells.Item(77, 6).Value = "set count"
$ws.Cells.Item(77, 7).Value = "DONE"
$ws.Cells.Item(77, 8).Value = 12
$ws.Cells.Item(77, 9).Value = 75.47
$ws.Cells.Item(77, 10).Value = 12
$ws.Cells.Item(77, 11).Value = 75.47
$ws.Cells.Item(77, 12).Value = 9.9
$ws.Cells.Item(77, 13).Value = 15.09
$ws.Cells.Item(77, 14).Value = 43847
$ws.Cells.Item(77, 14).NumberFormat = "yyyy/mm/dd"
$ws.Cells.Item(77, 15).Value = 0.86597222222222225
$ws.Cells.Item(77, 15).NumberFormat = "h:mm:ss"
$ws.Cells.Item(77, 16).Value = $false

Write-Host "A77 type:" $ws.Cells.Item(77,1).Value.GetType().FullName
